$d = $word.ActiveDocument

# Update the date paragraph heading
$d.Content.Find.Execute("2024-03-26 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-27 Wednesday", 2)

# Update the 25 multiplication answers in the table by directly targeting
# each cell (row, column) to avoid ambiguity from duplicate/overlapping
# text values that occur during the sequence of edits.
$t = $d.Tables.Item(1)

$newValues = @(
    @("76×76=5776", "74×48=3552", "79×12=948", "36×39=1404", "21×44=924"),
    @("41×81=3321", "39×33=1287", "20×22=440", "33×73=2409", "46×23=1058"),
    @("91×32=2912", "51×35=1785", "53×76=4028", "25×33=825", "39×36=1404"),
    @("61×13=793", "27×65=1755", "67×12=804", "18×97=1746", "87×55=4785"),
    @("94×91=8554", "34×98=3332", "27×93=2511", "17×27=459", "83×42=3486")
)

$dataRows = @(1, 5, 10, 15, 20)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $r = $dataRows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i][$c - 1]
    }
}
